$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Assign a value while preserving the cell as plain text, even when the
    # value looks like a number (e.g. "356.87"). We temporarily force a text
    # number format, assign the value, then restore the cell's original
    # style so no visible residue is left on the cell.
    $rng = $ws.Range($cellRef)
    $orig = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $orig
}

# Row 2 - Bitcoin
Set-TextValue "D2" "51.818.30"
Set-TextValue "E2" "  +0.22%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.820.55"
Set-TextValue "E3" "  +1.80%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.08%  "

# Row 5 - BNB
Set-TextValue "D5" "356.87"
Set-TextValue "E5" "  +3.54%  "

# Row 6 - Solana
Set-TextValue "D6" "111.87"
Set-TextValue "E6" "  -2.62%  "

# Row 7 - XRP
Set-TextValue "D7" "0.567"
Set-TextValue "E7" "  +3.52%  "

# Row 8 - USDC
Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  -0.04%  "

# Row 9 - Cardano
Set-TextValue "E9" "  +3.88%  "

# Row 10 - Avalanche
Set-TextValue "D10" "40.87"
Set-TextValue "E10" "  -3.91%  "

# Row 11 - Dogecoin
Set-TextValue "E11" "  +0.07%  "

# Row 12 - TRON
Set-TextValue "E12" "  +0.39%  "

# Row 13 - Chainlink
Set-TextValue "D13" "19.94"
Set-TextValue "E13" "  -0.26%  "

# Row 14 - Polkadot
Set-TextValue "E14" "  +1.55%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "3.261.10"
Set-TextValue "E15" "  +1.62%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "2.813.79"
Set-TextValue "E16" "  +0.78%  "

# Row 17 - Polygon
Set-TextValue "E17" "  +4.45%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "51.736.36"
Set-TextValue "E18" "  +0.13%  "

# Row 19 - Uniswap
Set-TextValue "D19" "7.56"
Set-TextValue "E19" "  +7.78%  "

# Row 20 - ImmutableX
Set-TextValue "D20" "3.13"
Set-TextValue "E20" "  -2.83%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue "D21" "13.36"
Set-TextValue "E21" "  +1.17%  "

# Row 22 - ShibaInu
Set-TextValue "D22" "0.0₃0990"
Set-TextValue "E22" "  +1.36%  "

# Row 23 - Litecoin
Set-TextValue "D23" "69.97"
Set-TextValue "E23" "  +0.07%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "268.03"
Set-TextValue "E24" "  -1.94%  "

# Row 25 - PancakeSwap
Set-TextValue "D25" "2.81"
Set-TextValue "E25" "  +1.51%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "26.98"
Set-TextValue "E26" "  +1.99%  "

# Row 27 - Dai
Set-TextValue "E27" "  +0.10%  "

# Row 28 - Cosmos
Set-TextValue "D28" "10.29"
Set-TextValue "E28" "  +0.87%  "

# Row 29 - Toncoin
Set-TextValue "D29" "2.25"
Set-TextValue "E29" "  +0.82%  "

# Row 30 - VeChain
Set-TextValue "D30" "0.0481"
Set-TextValue "E30" "  +26.47%  "

# Row 31 - Kaspa
Set-TextValue "D31" "0.140"
Set-TextValue "E31" "  -0.18%  "

# Row 32 - OKB
Set-TextValue "D32" "52.51"
Set-TextValue "E32" "  +5.09%  "

# Row 33 - InjectiveProtocol
Set-TextValue "D33" "34.55"
Set-TextValue "E33" "  +0.43%  "

# Row 34 - Filecoin
Set-TextValue "D34" "5.88"
Set-TextValue "E34" "  +3.21%  "

# Row 35 - RenderToken
Set-TextValue "E35" "  +8.94%  "

# Row 36 - Hedera
Set-TextValue "E36" "  +3.21%  "

# Row 37 - FirstDigitalUSD
Set-TextValue "D37" "0.999"
Set-TextValue "E37" "  -0.15%  "

# Row 38 - LidoDAOToken
Set-TextValue "D38" "3.31"
Set-TextValue "E38" "  +2.94%  "

# Row 39 - ARBITRUM
Set-TextValue "E39" "  -2.57%  "

# Row 40 - Celestia
Set-TextValue "D40" "18.35"
Set-TextValue "E40" "  -3.02%  "

# Row 41 - Stellar
Set-TextValue "E41" "  +0.50%  "

# Row 42 & 43 - EnergySwap / Stacks swap positions
Set-TextValue "B42" "Stacks"
Set-TextValue "C42" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D42" "2.53"
Set-TextValue "E42" "  -4.69%  "

Set-TextValue "B43" "EnergySwap"
Set-TextValue "C43" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D43" "23.15"
Set-TextValue "E43" "  -0.83%  "

# Row 44 - Monero
Set-TextValue "D44" "124.71"
Set-TextValue "E44" "  -2.10%  "

# Row 45 - WEMIXToken
Set-TextValue "E45" "  -2.68%  "

# Row 46 - Maker
Set-TextValue "D46" "2.090.69"
Set-TextValue "E46" "  +1.25%  "

# Row 47 - NEARProtocol
Set-TextValue "D47" "3.35"
Set-TextValue "E47" "  +1.49%  "

# Row 48 - ApeXProtocol
Set-TextValue "E48" "  +1.10%  "

# Row 49 - THORChain
Set-TextValue "D49" "5.97"
Set-TextValue "E49" "  +8.09%  "

# Row 50 - SEI
Set-TextValue "D50" "0.981"
Set-TextValue "E50" "  +10.99%  "

# Row 51 - FraxShare
Set-TextValue "D51" "9.04"
Set-TextValue "E51" "  +1.92%  "
